$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.401172666666667
$ws.Range("H2").Value = 10.203518
$ws.Range("I2").Value = 0.5101677883321656
$ws.Range("J2").Value = 0.5101677883321655
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 161.7750676666667
$ws.Range("N2").Value = 485.325203
$ws.Range("O2").Value = 0.9790864123038654
$ws.Range("P2").Value = 0.9790864123038654
$ws.Range("Q2").Value = 550.2249382960172
$ws.Range("R2").Value = 4952.024444664155
$ws.Range("S2").Value = 0.4994983495511378
$ws.Range("T2").Value = 0.4994983495511377
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.401172666666667
$ws.Range("H3").Value = 10.203518
$ws.Range("I3").Value = 0.5101677883321656
$ws.Range("J3").Value = 0.5101677883321655
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.67894
$ws.Range("N3").Value = 2.03682
$ws.Range("O3").Value = 0.004109044356199978
$ws.Range("P3").Value = 0.004109044356199979
$ws.Range("Q3").Value = 2.309192170306667
$ws.Range("R3").Value = 20.78272953276
$ws.Range("S3").Value = 0.00209630207136131
$ws.Range("T3").Value = 0.00209630207136131
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.401172666666667
$ws.Range("H4").Value = 10.203518
$ws.Range("I4").Value = 0.5101677883321656
$ws.Range("J4").Value = 0.5101677883321655
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.763201333333333
$ws.Range("N4").Value = 5.289604
$ws.Range("O4").Value = 0.01067115280816804
$ws.Range("P4").Value = 0.01067115280816804
$ws.Range("Q4").Value = 5.996952180763556
$ws.Range("R4").Value = 53.972569626872
$ws.Range("S4").Value = 0.005444078427097668
$ws.Range("T4").Value = 0.005444078427097667
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.401172666666667
$ws.Range("H5").Value = 10.203518
$ws.Range("I5").Value = 0.5101677883321656
$ws.Range("J5").Value = 0.5101677883321655
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.013424
$ws.Range("N5").Value = 3.040272
$ws.Range("O5").Value = 0.006133390531766587
$ws.Range("P5").Value = 0.006133390531766588
$ws.Range("Q5").Value = 3.446830008544
$ws.Range("R5").Value = 31.021470076896
$ws.Range("S5").Value = 0.003129058282568805
$ws.Range("T5").Value = 0.003129058282568804
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.579868
$ws.Range("H6").Value = 7.739604
$ws.Range("I6").Value = 0.3869740471126509
$ws.Range("J6").Value = 0.3869740471126508
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 161.7750676666667
$ws.Range("N6").Value = 485.325203
$ws.Range("O6").Value = 0.9790864123038654
$ws.Range("P6").Value = 0.9790864123038654
$ws.Range("Q6").Value = 417.358320271068
$ws.Range("R6").Value = 3756.224882439612
$ws.Range("S6").Value = 0.3788810314422323
$ws.Range("T6").Value = 0.3788810314422323
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.579868
$ws.Range("H7").Value = 7.739604
$ws.Range("I7").Value = 0.3869740471126509
$ws.Range("J7").Value = 0.3869740471126508
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.67894
$ws.Range("N7").Value = 2.03682
$ws.Range("O7").Value = 0.004109044356199978
$ws.Range("P7").Value = 0.004109044356199979
$ws.Range("Q7").Value = 1.75157557992
$ws.Range("R7").Value = 15.76418021928
$ws.Range("S7").Value = 0.001590093524284103
$ws.Range("T7").Value = 0.001590093524284103
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.579868
$ws.Range("H8").Value = 7.739604
$ws.Range("I8").Value = 0.3869740471126509
$ws.Range("J8").Value = 0.3869740471126508
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.763201333333333
$ws.Range("N8").Value = 5.289604
$ws.Range("O8").Value = 0.01067115280816804
$ws.Range("P8").Value = 0.01067115280816804
$ws.Range("Q8").Value = 4.548826697423999
$ws.Range("R8").Value = 40.939440276816
$ws.Range("S8").Value = 0.004129459189534316
$ws.Range("T8").Value = 0.004129459189534316
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.579868
$ws.Range("H9").Value = 7.739604
$ws.Range("I9").Value = 0.3869740471126509
$ws.Range("J9").Value = 0.3869740471126508
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.013424
$ws.Range("N9").Value = 3.040272
$ws.Range("O9").Value = 0.006133390531766587
$ws.Range("P9").Value = 0.006133390531766588
$ws.Range("Q9").Value = 2.614500148032
$ws.Range("R9").Value = 23.530501332288
$ws.Range("S9").Value = 0.00237346295660013
$ws.Range("T9").Value = 0.00237346295660013
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.6857320000000001
$ws.Range("H10").Value = 2.057196
$ws.Range("I10").Value = 0.1028581645551836
$ws.Range("J10").Value = 0.1028581645551836
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 161.7750676666667
$ws.Range("N10").Value = 485.325203
$ws.Range("O10").Value = 0.9790864123038654
$ws.Range("P10").Value = 0.9790864123038654
$ws.Range("Q10").Value = 110.9343407011987
$ws.Range("R10").Value = 998.4090663107881
$ws.Range("S10").Value = 0.1007070313104953
$ws.Range("T10").Value = 0.1007070313104953
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.6857320000000001
$ws.Range("H11").Value = 2.057196
$ws.Range("I11").Value = 0.1028581645551836
$ws.Range("J11").Value = 0.1028581645551836
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.67894
$ws.Range("N11").Value = 2.03682
$ws.Range("O11").Value = 0.004109044356199978
$ws.Range("P11").Value = 0.004109044356199979
$ws.Range("Q11").Value = 0.4655708840800001
$ws.Range("R11").Value = 4.190137956720001
$ws.Range("S11").Value = 0.0004226487605545658
$ws.Range("T11").Value = 0.0004226487605545658
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.6857320000000001
$ws.Range("H12").Value = 2.057196
$ws.Range("I12").Value = 0.1028581645551836
$ws.Range("J12").Value = 0.1028581645551836
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 1.763201333333333
$ws.Range("N12").Value = 5.289604
$ws.Range("O12").Value = 0.01067115280816804
$ws.Range("P12").Value = 0.01067115280816804
$ws.Range("Q12").Value = 1.209083576709334
$ws.Range("R12").Value = 10.881752190384
$ws.Range("S12").Value = 0.001097615191536058
$ws.Range("T12").Value = 0.001097615191536058
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.6857320000000001
$ws.Range("H13").Value = 2.057196
$ws.Range("I13").Value = 0.1028581645551836
$ws.Range("J13").Value = 0.1028581645551836
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.013424
$ws.Range("N13").Value = 3.040272
$ws.Range("O13").Value = 0.006133390531766587
$ws.Range("P13").Value = 0.006133390531766588
$ws.Range("Q13").Value = 0.694937266368
$ws.Range("R13").Value = 6.254435397312
$ws.Range("S13").Value = 0.0006308692925976526
$ws.Range("T13").Value = 0.0006308692925976525
